$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# Sheet "COORD": remove the "William Baril" row (row 5)
# ------------------------------------------------------------------
$wsCoord = $wb.Worksheets.Item("COORD")

# Remove the hyperlink that lives on B5 before clearing the cell
foreach ($h in $wsCoord.Hyperlinks) {
    $addr = $h.Range.Address()
    if ($addr -eq '$B$5') {
        $h.Delete() | Out-Null
    }
}

# Clear the whole staff row - it no longer exists
$wsCoord.Range("A5:C5").ClearContents() | Out-Null

# ------------------------------------------------------------------
# Sheet "COORD_DEP": update the "Mecano" role label and add the new
# "Guillaume Julien" staff member as row 6
# ------------------------------------------------------------------
$wsCoordDep = $wb.Worksheets.Item("COORD_DEP")

# New row 6 (name + role entered first)
$wsCoordDep.Range("A6").Value = "Guillaume Julien"
$wsCoordDep.Range("D6").Value = "Mécano 3"

# Row 2's role note changes from "Mecano" to "Mecano 3"
$wsCoordDep.Range("D2").Value = "Mecano 3"

# Remaining new row 6 details
$wsCoordDep.Range("B6").Value = "gjulien18@hotmail.com"
$wsCoordDep.Range("C6").Value = "819-860-4596"

# C6 did not previously hold any data - copy the formatting used by
# the cells above it (column C / style used on rows 3-5)
$wsCoordDep.Range("C5").Copy() | Out-Null
$wsCoordDep.Range("C6").PasteSpecial(-4122) | Out-Null

# Add the mailto hyperlink on B6, then restore its normal formatting
# (adding a hyperlink re-styles the cell as a hyperlink by default)
$wsCoordDep.Hyperlinks.Add($wsCoordDep.Range("B6"), "mailto:gjulien18@hotmail.com") | Out-Null
$wsCoordDep.Range("B5").Copy() | Out-Null
$wsCoordDep.Range("B6").PasteSpecial(-4122) | Out-Null

$wsCoordDep.Application.CutCopyMode = $false

# ------------------------------------------------------------------
# Sheet "HORAIRE_DEP": replace the "TBD" placeholders with the new
# staff member's name
# ------------------------------------------------------------------
$wsHoraireDep = $wb.Worksheets.Item("HORAIRE_DEP")
$wsHoraireDep.Range("F2").Value = "Guillaume Julien"
$wsHoraireDep.Range("F5").Value = "Guillaume Julien"
$wsHoraireDep.Range("F6").Value = "Guillaume Julien"

# ------------------------------------------------------------------
# Restore the selection / active cell on each sheet, finishing on
# COORD so it stays the active tab
# ------------------------------------------------------------------
$wsCoordDep.Range("A8").Select() | Out-Null
$wsHoraireDep.Range("F14").Select() | Out-Null
$wsCoord.Range("A30").Select() | Out-Null
